$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -20.485
$ws.Range("A7").Value = -19.657
$ws.Range("C7").Value = -12.748
$ws.Range("C15").Value = -13.173
$ws.Range("A16").Value = -21.861
$ws.Range("D16").Value = -8.598000000000001
$ws.Range("D19").Value = -7.855
$ws.Range("C21").Value = -12.132
$ws.Range("C22").Value = -12.867
$ws.Range("C23").Value = -12.141
$ws.Range("A28").Value = -21.86
$ws.Range("A29").Value = -21.246
$ws.Range("A32").Value = -21.87
$ws.Range("C34").Value = -12.089
$ws.Range("E34").Value = 17.142
$ws.Range("D36").Value = -7.827
$ws.Range("A40").Value = -19.912
$ws.Range("C43").Value = -12.759
$ws.Range("E43").Value = 17.074
$ws.Range("C45").Value = -13.03
$ws.Range("D46").Value = -8.372
$ws.Range("E48").Value = 16.938
$ws.Range("C50").Value = -13.018
$ws.Range("D50").Value = -8.097
$ws.Range("C51").Value = -10.866
$ws.Range("A52").Value = -21.94
$ws.Range("A57").Value = -22.068
$ws.Range("A66").Value = -21.738
$ws.Range("C66").Value = -10.916
$ws.Range("C67").Value = -11.286
$ws.Range("E70").Value = 17.651
$ws.Range("E73").Value = 16.661
$ws.Range("C79").Value = -11.871
$ws.Range("C84").Value = -14.098
$ws.Range("E87").Value = 16.509
$ws.Range("C92").Value = -11.434
$ws.Range("E92").Value = 17.555
$ws.Range("D95").Value = -7.536000000000001
$ws.Range("C97").Value = -12.619
$ws.Range("D97").Value = -8.395999999999999
$ws.Range("A100").Value = -22.067
$ws.Range("E101").Value = 16.673
